# Weekly update: insert a new record at row 541 (Feria Lagunitas de Puerto
# Montt - Lechuga, Escarola) and push the existing rows 541-619 down to
# 542-620, matching the new dataset dump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 541; this shifts rows 541:619
# down to 542:620 and extends the used range to A1:R620.
$ws.Rows("541:541").Insert()

# Populate the newly inserted row 541 with the new record.
$ws.Range("A541").Value = 4
$ws.Range("B541").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C541").Value = "Los Lagos"
$ws.Range("D541").Value = 44776
$ws.Range("E541").Value = 10
$ws.Range("F541").Value = 100112033
$ws.Range("G541").Value = "Lechuga"
$ws.Range("H541").Value = "Escarola"
$ws.Range("I541").Value = "Segunda"
$ws.Range("J541").Value = 120
$ws.Range("K541").Value = 12000
$ws.Range("L541").Value = 12000
$ws.Range("M541").Value = 12000
$ws.Range("N541").Value = '$/caja 18 unidades'
$ws.Range("O541").Value = "Región de Coquimbo"
$ws.Range("P541").Value = 667
$ws.Range("Q541").Value = 18
$ws.Range("R541").Value = "Hortaliza"
